# Update PLC data 2025-10-13 13:42:18
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 256
$ws.Range("C3").Value = 155145
$ws.Range("C4").Value = 146269
$ws.Range("C7").Value = 5.72
$ws.Range("C8").Value = 63.48
